$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.348.33"
$ws.Range("E2").Value = "  +2.24%  "
$ws.Range("D3").Value = "2.427.03"
$ws.Range("E3").Value = "  +3.29%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "556.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.05%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.536"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.76%  "
$ws.Range("D9").Value = "2.427.15"
$ws.Range("E9").Value = "  +3.39%  "
$ws.Range("E10").Value = "  +5.10%  "
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("E12").Value = "  +1.31%  "
$ws.Range("E13").Value = "  +3.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.65%  "
$ws.Range("E15").Value = "  +9.50%  "
$ws.Range("D16").Value = "2.866.29"
$ws.Range("E16").Value = "  +3.32%  "
$ws.Range("D17").Value = "62.233.68"
$ws.Range("E17").Value = "  +2.41%  "
$ws.Range("D18").Value = "2.424.93"
$ws.Range("E18").Value = "  +3.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.10"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "325.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.09%  "
$ws.Range("E21").Value = "  +1.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.76"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.37%  "
$ws.Range("E23").Value = "  +1.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.86%  "
$ws.Range("E26").Value = "  +10.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "570.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +14.44%  "
$ws.Range("D28").Value = "2.541.74"
$ws.Range("E28").Value = "  +3.30%  "
$ws.Range("D29").Value = "0.0₃0946"
$ws.Range("E29").Value = "  +9.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.42"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.69%  "
$ws.Range("E32").Value = "  +5.29%  "
$ws.Range("E33").Value = "  +0.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.87"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.41%  "
$ws.Range("E35").Value = "  +4.85%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.73"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.71%  "
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.82"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.385"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.05%  "
$ws.Range("E40").Value = "  +4.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.86"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.96%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "147.71"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.26%  "
$ws.Range("E43").Value = "  +0.43%  "
$ws.Range("E44").Value = "  +3.19%  "
$ws.Range("E45").Value = "  +12.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "152.13"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.64"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0545"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.44"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.591"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.68%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0915"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.36%  "
